# Updated cryptos list with latest price/volume data (and two coin-rank swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '63.050.52' },
    @{ Cell = 'E2'; Value = '  -0.66%  ' },
    @{ Cell = 'D3'; Value = '2.552.62' },
    @{ Cell = 'E3'; Value = '  -0.08%  ' },
    @{ Cell = 'E4'; Value = '  -0.01%  ' },
    @{ Cell = 'D5'; Value = '581.40' },
    @{ Cell = 'E5'; Value = '  +1.65%  ' },
    @{ Cell = 'D6'; Value = '146.84' },
    @{ Cell = 'E6'; Value = '  -2.62%  ' },
    @{ Cell = 'E7'; Value = '  -0.02%  ' },
    @{ Cell = 'D8'; Value = '0.584' },
    @{ Cell = 'E8'; Value = '  -0.52%  ' },
    @{ Cell = 'E9'; Value = '  -0.24%  ' },
    @{ Cell = 'D10'; Value = '5.57' },
    @{ Cell = 'E10'; Value = '  -2.06%  ' },
    @{ Cell = 'E11'; Value = '  -0.17%  ' },
    @{ Cell = 'D12'; Value = '0.355' },
    @{ Cell = 'E12'; Value = '  -0.85%  ' },
    @{ Cell = 'D13'; Value = '27.56' },
    @{ Cell = 'E13'; Value = '  -2.18%  ' },
    @{ Cell = 'D14'; Value = '3.009.89' },
    @{ Cell = 'E14'; Value = '  -0.05%  ' },
    @{ Cell = 'D15'; Value = '62.969.54' },
    @{ Cell = 'E15'; Value = '  -0.61%  ' },
    @{ Cell = 'D16'; Value = '0.0000144' },
    @{ Cell = 'E16'; Value = '  +0.14%  ' },
    @{ Cell = 'D17'; Value = '2.550.12' },
    @{ Cell = 'E17'; Value = '  -0.84%  ' },
    @{ Cell = 'D18'; Value = '11.35' },
    @{ Cell = 'E18'; Value = '  -2.31%  ' },
    @{ Cell = 'D19'; Value = '339.35' },
    @{ Cell = 'E19'; Value = '  -0.29%  ' },
    @{ Cell = 'D20'; Value = '4.35' },
    @{ Cell = 'E20'; Value = '  +0.43%  ' },
    @{ Cell = 'E21'; Value = '  -0.59%  ' },
    @{ Cell = 'D22'; Value = '1.00' },
    @{ Cell = 'E22'; Value = '  -0.08%  ' },
    @{ Cell = 'D23'; Value = '65.68' },
    @{ Cell = 'E23'; Value = '  -0.44%  ' },
    @{ Cell = 'D24'; Value = '2.676.79' },
    @{ Cell = 'E24'; Value = '  +0.34%  ' },
    @{ Cell = 'B25'; Value = 'Fetch.AI' },
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' },
    @{ Cell = 'D25'; Value = '1.63' },
    @{ Cell = 'E25'; Value = '  +0.98%  ' },
    @{ Cell = 'B26'; Value = 'Kaspa' },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' },
    @{ Cell = 'D26'; Value = '0.169' },
    @{ Cell = 'E26'; Value = '  -0.17%  ' },
    @{ Cell = 'E27'; Value = '  -4.12%  ' },
    @{ Cell = 'D28'; Value = '1.00' },
    @{ Cell = 'E28'; Value = '  +0.14%  ' },
    @{ Cell = 'D29'; Value = '8.34' },
    @{ Cell = 'E29'; Value = '  -1.77%  ' },
    @{ Cell = 'D30'; Value = '7.77' },
    @{ Cell = 'E30'; Value = '  +5.66%  ' },
    @{ Cell = 'E31'; Value = '  +4.75%  ' },
    @{ Cell = 'D32'; Value = '0.0₃0816' },
    @{ Cell = 'E32'; Value = '  -0.51%  ' },
    @{ Cell = 'D33'; Value = '177.54' },
    @{ Cell = 'E33'; Value = '  -0.32%  ' },
    @{ Cell = 'B34'; Value = 'Bittensor' },
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' },
    @{ Cell = 'D34'; Value = '428.07' },
    @{ Cell = 'E34'; Value = '  +0.85%  ' },
    @{ Cell = 'B35'; Value = 'ImmutableX' },
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' },
    @{ Cell = 'D35'; Value = '1.57' },
    @{ Cell = 'E35'; Value = '  -1.46%  ' },
    @{ Cell = 'B36'; Value = 'EthereumClassic' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' },
    @{ Cell = 'D36'; Value = '19.18' },
    @{ Cell = 'E36'; Value = '  +0.87%  ' },
    @{ Cell = 'B37'; Value = 'PolygonEcosystemToken' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol' },
    @{ Cell = 'D37'; Value = '0.401' },
    @{ Cell = 'E37'; Value = '  -0.61%  ' },
    @{ Cell = 'D38'; Value = '4.39' },
    @{ Cell = 'E38'; Value = '  -1.86%  ' },
    @{ Cell = 'E39'; Value = '  +0.02%  ' },
    @{ Cell = 'E40'; Value = '  -1.90%  ' },
    @{ Cell = 'E41'; Value = '  +0.00%  ' },
    @{ Cell = 'D42'; Value = '39.62' },
    @{ Cell = 'E42'; Value = '  +0.19%  ' },
    @{ Cell = 'D43'; Value = '150.91' },
    @{ Cell = 'E43'; Value = '  -1.98%  ' },
    @{ Cell = 'D44'; Value = '3.77' },
    @{ Cell = 'E44'; Value = '  -0.44%  ' },
    @{ Cell = 'D45'; Value = '20.79' },
    @{ Cell = 'E45'; Value = '  -0.62%  ' },
    @{ Cell = 'D46'; Value = '0.0542' },
    @{ Cell = 'E46'; Value = '  +3.21%  ' },
    @{ Cell = 'B47'; Value = 'Stellar' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Cell = 'D47'; Value = '0.0973' },
    @{ Cell = 'E47'; Value = '  +0.51%  ' },
    @{ Cell = 'B48'; Value = 'Mantle' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt' },
    @{ Cell = 'D48'; Value = '0.601' },
    @{ Cell = 'E48'; Value = '  -1.58%  ' },
    @{ Cell = 'D49'; Value = '0.0240' },
    @{ Cell = 'E49'; Value = '  +0.28%  ' },
    @{ Cell = 'D50'; Value = '18.27' },
    @{ Cell = 'E50'; Value = '  -1.71%  ' },
    @{ Cell = 'E51'; Value = '  -0.43%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = '@'   # preserve text formatting (e.g. "1.00", leading/trailing spaces, dotted thousands)
    $cell.Value = $u.Value
}
